$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.100.10"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "2.444.81"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "2.440.21"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("E11").Value = "  +2.95%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("E15").Value = "  +1.31%  "

$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").Value = "62.102.53"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").Value = "2.436.98"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("E22").Value = "  -0.68%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  -5.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "593.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.29%  "

$ws.Range("D28").Value = "0.0₃0970"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").Value = "2.566.76"
$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("E30").Value = "  +2.98%  "

$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("E33").Value = "  +1.71%  "

$ws.Range("E34").Value = "  +0.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("E43").Value = "  -0.94%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "141.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").Value = "0.0₆0269"
$ws.Range("E48").Value = "  +21.17%  "

$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "

# Rows 37-39 reorder: Monero, ImmutableX, PolygonEcosystemToken -> ImmutableX, PolygonEcosystemToken, Monero
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.67%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.374"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.34%  "
